$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 15222.286  # H70: 10815.2 -> 15222.286
$ws.Cells.Item(70, 10).Value = 21031.2  # J70: 13344 -> 21031.2
$ws.Cells.Item(70, 12).Value = 63093.60000000001  # L70: 40032 -> 63093.60000000001
$ws.Cells.Item(70, 14).Value = -63633.60000000001  # N70: -40572 -> -63633.60000000001
$ws.Cells.Item(73, 8).Value = 15222.286  # H73: 10815.2 -> 15222.286
$ws.Cells.Item(73, 10).Value = 21031.2  # J73: 13344 -> 21031.2
$ws.Cells.Item(73, 12).Value = 63093.60000000001  # L73: 40032 -> 63093.60000000001
$ws.Cells.Item(73, 14).Value = -64965.60000000001  # N73: -41904 -> -64965.60000000001
$ws.Cells.Item(100, 8).Value = 1556.1818  # H100: 1471.9 -> 1556.1818
$ws.Cells.Item(100, 9).Value = 1411.8  # I100: 1471.9 -> 1411.8
$ws.Cells.Item(100, 10).Value = 3000  # J100: 0 -> 3000
$ws.Cells.Item(100, 11).Value = 1411.8  # K100: 1471.9 -> 1411.8
$ws.Cells.Item(100, 12).Value = 3000  # L100: 0 -> 3000
$ws.Cells.Item(100, 13).Value = -870.8  # M100: -930.9000000000001 -> -870.8
$ws.Cells.Item(100, 14).Value = -4082  # N100: None -> -4082
$ws.Cells.Item(107, 8).Value = 1228.6  # H107: 1311 -> 1228.6
$ws.Cells.Item(107, 9).Value = 949.5  # I107: 966.3333 -> 949.5
$ws.Cells.Item(107, 11).Value = 949.5  # K107: 966.3333 -> 949.5
$ws.Cells.Item(107, 13).Value = 970.5  # M107: 953.6667 -> 970.5
$ws.Cells.Item(121, 8).Value = 999  # H121: 0 -> 999
$ws.Cells.Item(121, 10).Value = 999  # J121: 0 -> 999
$ws.Cells.Item(121, 12).Value = 2997  # L121: 0 -> 2997
$ws.Cells.Item(121, 14).Value = -6491  # N121: None -> -6491
$ws.Cells.Item(132, 8).Value = 1166  # H132: 1177.7742 -> 1166
$ws.Cells.Item(132, 9).Value = 1114.2413  # I132: 1137.3334 -> 1114.2413
$ws.Cells.Item(132, 10).Value = 1666.3334  # J132: 1450.75 -> 1666.3334
$ws.Cells.Item(132, 11).Value = 3342.7239  # K132: 3412.0002 -> 3342.7239
$ws.Cells.Item(132, 12).Value = 4999.0002  # L132: 4352.25 -> 4999.0002
$ws.Cells.Item(132, 13).Value = -812.7239  # M132: -882.0001999999999 -> -812.7239
$ws.Cells.Item(132, 14).Value = -10059.0002  # N132: -9412.25 -> -10059.0002
$ws.Cells.Item(138, 8).Value = 2925.8462  # H138: 2991.08 -> 2925.8462
$ws.Cells.Item(138, 9).Value = 3039.3333  # I138: 3197.9092 -> 3039.3333
$ws.Cells.Item(138, 11).Value = 9117.999899999999  # K138: 9593.7276 -> 9117.999899999999
$ws.Cells.Item(138, 13).Value = -3977.999899999999  # M138: -4453.7276 -> -3977.999899999999
$ws.Cells.Item(139, 8).Value = 48550  # H139: 48580 -> 48550
$ws.Cells.Item(139, 10).Value = 48550  # J139: 48580 -> 48550
$ws.Cells.Item(139, 12).Value = 48550  # L139: 48580 -> 48550
$ws.Cells.Item(139, 14).Value = -58830  # N139: -58860 -> -58830
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4649.0713  # H32: 4973.154 -> 4649.0713
$ws.Cells.Item(32, 9).Value = 3325.7942  # I32: 3605.4517 -> 3325.7942
$ws.Cells.Item(32, 11).Value = 3325.7942  # K32: 3605.4517 -> 3325.7942
$ws.Cells.Item(32, 13).Value = -3038.7942  # M32: -3318.4517 -> -3038.7942
$ws.Cells.Item(44, 8).Value = 30000  # H44: 0 -> 30000
$ws.Cells.Item(44, 10).Value = 30000  # J44: 0 -> 30000
$ws.Cells.Item(44, 12).Value = 30000  # L44: 0 -> 30000
$ws.Cells.Item(44, 14).Value = -30976  # N44: None -> -30976
$ws.Cells.Item(61, 8).Value = 5159.4614  # H61: 5455.4614 -> 5159.4614
$ws.Cells.Item(61, 9).Value = 3107.4  # I61: 3158 -> 3107.4
$ws.Cells.Item(61, 10).Value = 11999.667  # J61: 10624.75 -> 11999.667
$ws.Cells.Item(61, 11).Value = 3107.4  # K61: 3158 -> 3107.4
$ws.Cells.Item(61, 12).Value = 11999.667  # L61: 10624.75 -> 11999.667
$ws.Cells.Item(61, 13).Value = -2895.4  # M61: -2946 -> -2895.4
$ws.Cells.Item(61, 14).Value = -12423.667  # N61: -11048.75 -> -12423.667
$ws.Cells.Item(110, 8).Value = 2563.4546  # H110: 2454.3635 -> 2563.4546
$ws.Cells.Item(110, 9).Value = 1381.7142  # I110: 1434 -> 1381.7142
$ws.Cells.Item(110, 10).Value = 4631.5  # J110: 5175.3335 -> 4631.5
$ws.Cells.Item(110, 11).Value = 1381.7142  # K110: 1434 -> 1381.7142
$ws.Cells.Item(110, 12).Value = 4631.5  # L110: 5175.3335 -> 4631.5
$ws.Cells.Item(110, 13).Value = 663.2858000000001  # M110: 611 -> 663.2858000000001
$ws.Cells.Item(110, 14).Value = -8721.5  # N110: -9265.333500000001 -> -8721.5
$ws.Cells.Item(132, 8).Value = 1757.1666  # H132: 1602.9642 -> 1757.1666
$ws.Cells.Item(132, 9).Value = 1061.6428  # I132: 976.3333 -> 1061.6428
$ws.Cells.Item(132, 11).Value = 3184.9284  # K132: 2928.9999 -> 3184.9284
$ws.Cells.Item(132, 13).Value = -654.9284000000002  # M132: -398.9998999999998 -> -654.9284000000002
$ws.Cells.Item(136, 8).Value = 5159.4614  # H136: 5455.4614 -> 5159.4614
$ws.Cells.Item(136, 9).Value = 3107.4  # I136: 3158 -> 3107.4
$ws.Cells.Item(136, 10).Value = 11999.667  # J136: 10624.75 -> 11999.667
$ws.Cells.Item(136, 11).Value = 9322.200000000001  # K136: 9474 -> 9322.200000000001
$ws.Cells.Item(136, 12).Value = 35999.001  # L136: 31874.25 -> 35999.001
$ws.Cells.Item(136, 13).Value = -6772.200000000001  # M136: -6924 -> -6772.200000000001
$ws.Cells.Item(136, 14).Value = -41099.001  # N136: -36974.25 -> -41099.001
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(28, 8).Value = 0  # H28: 30542 -> 0
$ws.Cells.Item(28, 10).Value = 0  # J28: 30542 -> 0
$ws.Cells.Item(28, 12).Value = 0  # L28: 30542 -> 0
$ws.Cells.Item(28, 14).ClearContents()  # N28 was -31130
$ws.Cells.Item(134, 8).Value = 19868.5  # H134: 17986.055 -> 19868.5
$ws.Cells.Item(134, 9).Value = 21816  # I134: 21196.363 -> 21816
$ws.Cells.Item(134, 10).Value = 14999.75  # J134: 12941.286 -> 14999.75
$ws.Cells.Item(134, 11).Value = 65448  # K134: 63589.08900000001 -> 65448
$ws.Cells.Item(134, 12).Value = 44999.25  # L134: 38823.858 -> 44999.25
$ws.Cells.Item(134, 13).Value = -62913  # M134: -61054.08900000001 -> -62913
$ws.Cells.Item(134, 14).Value = -50069.25  # N134: -43893.858 -> -50069.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2772.8538  # H31: 2571.2273 -> 2772.8538
$ws.Cells.Item(31, 9).Value = 1903.2  # I31: 1710.8334 -> 1903.2
$ws.Cells.Item(31, 10).Value = 4131.6875  # J31: 4414.9287 -> 4131.6875
$ws.Cells.Item(31, 11).Value = 1903.2  # K31: 1710.8334 -> 1903.2
$ws.Cells.Item(31, 12).Value = 4131.6875  # L31: 4414.9287 -> 4131.6875
$ws.Cells.Item(31, 13).Value = -1608.2  # M31: -1415.8334 -> -1608.2
$ws.Cells.Item(31, 14).Value = -4721.6875  # N31: -5004.9287 -> -4721.6875
$ws.Cells.Item(34, 8).Value = 2772.8538  # H34: 2571.2273 -> 2772.8538
$ws.Cells.Item(34, 9).Value = 1903.2  # I34: 1710.8334 -> 1903.2
$ws.Cells.Item(34, 10).Value = 4131.6875  # J34: 4414.9287 -> 4131.6875
$ws.Cells.Item(34, 11).Value = 1903.2  # K34: 1710.8334 -> 1903.2
$ws.Cells.Item(34, 12).Value = 4131.6875  # L34: 4414.9287 -> 4131.6875
$ws.Cells.Item(34, 13).Value = -1701.2  # M34: -1508.8334 -> -1701.2
$ws.Cells.Item(34, 14).Value = -4535.6875  # N34: -4818.9287 -> -4535.6875
$ws.Cells.Item(62, 8).Value = 2565  # H62: 2410.6365 -> 2565
$ws.Cells.Item(62, 9).Value = 2374  # I62: 2311.6667 -> 2374
$ws.Cells.Item(62, 10).Value = 2883.3333  # J62: 2529.4 -> 2883.3333
$ws.Cells.Item(62, 11).Value = 2374  # K62: 2311.6667 -> 2374
$ws.Cells.Item(62, 12).Value = 2883.3333  # L62: 2529.4 -> 2883.3333
$ws.Cells.Item(62, 13).Value = -1750  # M62: -1687.6667 -> -1750
$ws.Cells.Item(62, 14).Value = -4131.3333  # N62: -3777.4 -> -4131.3333
$ws.Cells.Item(65, 8).Value = 2565  # H65: 2410.6365 -> 2565
$ws.Cells.Item(65, 9).Value = 2374  # I65: 2311.6667 -> 2374
$ws.Cells.Item(65, 10).Value = 2883.3333  # J65: 2529.4 -> 2883.3333
$ws.Cells.Item(65, 11).Value = 11870  # K65: 11558.3335 -> 11870
$ws.Cells.Item(65, 12).Value = 14416.6665  # L65: 12647 -> 14416.6665
$ws.Cells.Item(65, 13).Value = -8750  # M65: -8438.333500000001 -> -8750
$ws.Cells.Item(65, 14).Value = -20656.6665  # N65: -18887 -> -20656.6665
$ws.Cells.Item(107, 8).Value = 459.05884  # H107: 361.73914 -> 459.05884
$ws.Cells.Item(107, 9).Value = 400.3846  # I107: 301.10526 -> 400.3846
$ws.Cells.Item(107, 11).Value = 400.3846  # K107: 301.10526 -> 400.3846
$ws.Cells.Item(107, 13).Value = 1519.6154  # M107: 1618.89474 -> 1519.6154
$ws.Cells.Item(122, 8).Value = 1951  # H122: 2046.3846 -> 1951
$ws.Cells.Item(122, 9).Value = 1884.5  # I122: 1991.1818 -> 1884.5
$ws.Cells.Item(122, 11).Value = 5653.5  # K122: 5973.5454 -> 5653.5
$ws.Cells.Item(122, 13).Value = -3203.5  # M122: -3523.5454 -> -3203.5
$ws.Cells.Item(134, 8).Value = 743.76666  # H134: 764.76666 -> 743.76666
$ws.Cells.Item(134, 9).Value = 743.76666  # I134: 764.76666 -> 743.76666
$ws.Cells.Item(134, 11).Value = 2231.29998  # K134: 2294.29998 -> 2231.29998
$ws.Cells.Item(134, 13).Value = 303.7000200000002  # M134: 240.7000200000002 -> 303.7000200000002
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2408063  # H132: 2408213 -> 2408063
$ws.Cells.Item(132, 9).Value = 2961693  # I132: 2961877.8 -> 2961693
$ws.Cells.Item(132, 11).Value = 8885079  # K132: 8885633.399999999 -> 8885079
$ws.Cells.Item(132, 13).Value = -8882549  # M132: -8883103.399999999 -> -8882549
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 2471.9333  # H68: 2434.2856 -> 2471.9333
$ws.Cells.Item(68, 10).Value = 3666.3333  # J68: 4000 -> 3666.3333
$ws.Cells.Item(68, 12).Value = 3666.3333  # L68: 4000 -> 3666.3333
$ws.Cells.Item(68, 14).Value = -5164.3333  # N68: -5498 -> -5164.3333
$ws.Cells.Item(71, 8).Value = 2471.9333  # H71: 2434.2856 -> 2471.9333
$ws.Cells.Item(71, 10).Value = 3666.3333  # J71: 4000 -> 3666.3333
$ws.Cells.Item(71, 12).Value = 18331.6665  # L71: 20000 -> 18331.6665
$ws.Cells.Item(71, 14).Value = -25819.6665  # N71: -27488 -> -25819.6665
$ws.Cells.Item(93, 8).Value = 547.3  # H93: 376.09525 -> 547.3
$ws.Cells.Item(93, 9).Value = 467.57144  # I93: 332.84616 -> 467.57144
$ws.Cells.Item(93, 10).Value = 733.3333  # J93: 446.375 -> 733.3333
$ws.Cells.Item(93, 11).Value = 467.57144  # K93: 332.84616 -> 467.57144
$ws.Cells.Item(93, 12).Value = 733.3333  # L93: 446.375 -> 733.3333
$ws.Cells.Item(93, 13).Value = 780.4285600000001  # M93: 915.1538399999999 -> 780.4285600000001
$ws.Cells.Item(93, 14).Value = -3229.3333  # N93: -2942.375 -> -3229.3333
$ws.Cells.Item(122, 8).Value = 12666.833  # H122: 11429 -> 12666.833
$ws.Cells.Item(122, 9).Value = 11500.25  # I122: 10000.6 -> 11500.25
$ws.Cells.Item(122, 11).Value = 34500.75  # K122: 30001.8 -> 34500.75
$ws.Cells.Item(122, 13).Value = -32050.75  # M122: -27551.8 -> -32050.75
$ws.Cells.Item(132, 8).Value = 3262.7778  # H132: 3495.625 -> 3262.7778
$ws.Cells.Item(132, 9).Value = 2632.8333  # I132: 2879.4 -> 2632.8333
$ws.Cells.Item(132, 11).Value = 7898.499899999999  # K132: 8638.200000000001 -> 7898.499899999999
$ws.Cells.Item(132, 13).Value = -5368.499899999999  # M132: -6108.200000000001 -> -5368.499899999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(48, 8).Value = 68000  # H48: 0 -> 68000
$ws.Cells.Item(48, 10).Value = 68000  # J48: 0 -> 68000
$ws.Cells.Item(48, 12).Value = 68000  # L48: 0 -> 68000
$ws.Cells.Item(48, 14).Value = -69138  # N48: None -> -69138
$ws.Cells.Item(70, 8).Value = 29400  # H70: 29039.4 -> 29400
$ws.Cells.Item(70, 10).Value = 29400  # J70: 29039.4 -> 29400
$ws.Cells.Item(70, 12).Value = 29400  # L70: 29039.4 -> 29400
$ws.Cells.Item(70, 14).Value = -30030  # N70: -29669.4 -> -30030
$ws.Cells.Item(73, 8).Value = 29400  # H73: 29039.4 -> 29400
$ws.Cells.Item(73, 10).Value = 29400  # J73: 29039.4 -> 29400
$ws.Cells.Item(73, 12).Value = 29400  # L73: 29039.4 -> 29400
$ws.Cells.Item(73, 14).Value = -31584  # N73: -31223.4 -> -31584
$ws.Cells.Item(122, 8).Value = 61075.23  # H122: 61075.31 -> 61075.23
$ws.Cells.Item(122, 9).Value = 71788.91  # I122: 71789 -> 71788.91
$ws.Cells.Item(122, 11).Value = 215366.73  # K122: 215367 -> 215366.73
$ws.Cells.Item(122, 13).Value = -212916.73  # M122: -212917 -> -212916.73
$ws.Cells.Item(132, 8).Value = 3581.1667  # H132: 3776.4443 -> 3581.1667
$ws.Cells.Item(132, 9).Value = 2999.5  # I132: 2999.4 -> 2999.5
$ws.Cells.Item(132, 10).Value = 4162.8335  # J132: 4747.75 -> 4162.8335
$ws.Cells.Item(132, 11).Value = 8998.5  # K132: 8998.200000000001 -> 8998.5
$ws.Cells.Item(132, 12).Value = 12488.5005  # L132: 14243.25 -> 12488.5005
$ws.Cells.Item(132, 13).Value = -6468.5  # M132: -6468.200000000001 -> -6468.5
$ws.Cells.Item(132, 14).Value = -17548.5005  # N132: -19303.25 -> -17548.5005
